$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$ws.Range("C2").Value = 45207
$ws.Range("C3").Value = 45207
$ws.Range("C4").Value = 45207
$ws.Range("C5").Value = 45207
